$d = $word.ActiveDocument

# Locate the paragraph containing the MATLAB "fplot(x,z)" line so we can
# mark the new edit location with the _GoBack bookmark once the text is
# fixed (Word relocates _GoBack to the spot of the most recent edit).
$target = $d.Content
$target.Find.Execute("fplot(x,z)", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$editStart = $target.Start

# Fix the function name so it matches the actual MATLAB built-in "plot"
# (the previous "fplot(x,z)" was a typo referencing the wrong function).
$d.Content.Find.Execute("fplot(x,z)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "plot(x,z)", 2)

# Move the _GoBack bookmark (tracks the location of the last edit) from the
# old position to the start of the line that was just corrected.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$newRange = $d.Range($editStart, $editStart)
$d.Bookmarks.Add("_GoBack", $newRange)
